# Fix Training Data Issue:
# The "Date" column (BF) held the spreadsheet's own file-name-derived label
# ("6-29-2012-13") instead of an ISO date. Correct it to "2013-06-29" for
# every data row (BF2:BF31), same as the upstream fix.
#
# NumberFormat is forced to Text ("@") before the assignment so this engine
# doesn't auto-convert the ISO-looking string into a date serial, then
# ClearFormats() drops that temporary number-format tweak again so the
# cell's style/format is left exactly as it was (General / default style) -
# only the stored text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    $cell.NumberFormat = "@"
    $cell.Value = "2013-06-29"
    $cell.ClearFormats()
}

$wb.Save()
